$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 129, shifting rows 129:133 down to 130:134.
$ws.Rows.Item(129).Insert()

# Populate the new row 129 with the new weekly price entry (carries the
# same Mercado/Region/Producto metadata as the entry that used to sit
# at row 129, but with a new date, quality, volume and price figures).
$ws.Cells.Item(129, 1).Value = 10
$ws.Cells.Item(129, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(129, 3).Value = "La Araucanía"
$ws.Cells.Item(129, 4).Value = 44746
$ws.Cells.Item(129, 5).Value = 9
$ws.Cells.Item(129, 6).Value = "Fruta"
$ws.Cells.Item(129, 7).Value = 100104
$ws.Cells.Item(129, 8).Value = "Frutos de pepita"
$ws.Cells.Item(129, 9).Value = 100104001
$ws.Cells.Item(129, 10).Value = "Granada"
$ws.Cells.Item(129, 11).Value = "Wonderfull"
$ws.Cells.Item(129, 12).Value = "Especial"
$ws.Cells.Item(129, 13).Value = 200
$ws.Cells.Item(129, 14).Value = 13000
$ws.Cells.Item(129, 15).Value = 13000
$ws.Cells.Item(129, 16).Value = 13000
$ws.Cells.Item(129, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(129, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(129, 19).Value = 1300
$ws.Cells.Item(129, 20).Value = 10
